# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Rows 16-59 of sheet "Hoja1" list one "Periodo Mora" (column E, a
# YYMM text code) and its "Valor Mora" (column F) per row. The periods
# used to be listed newest-first (2002 down to 1607); this update
# re-sorts them oldest-first (1607 up to 2002) and refreshes the mora
# amount that goes with each period.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New chronological (oldest -> newest) period codes for rows 16..59.
$periods = @(
    "1607","1608","1609","1610","1611","1612",
    "1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
    "2001","2002"
)

# Updated "Valor Mora" amounts that line up with the periods above.
$valores = @(
    27560,27560,27560,27560,27560,27560,
    27560,27560,27560,27560,27560,27560,27560,27560,27560,27560,27560,27560,
    27560,27560,27560,27560,27560,27560,27560,27560,
    31249,31249,31249,31249,
    31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,
    31249,17708
)

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("E$row").Value = $periods[$i]
    $ws.Range("F$row").Value = $valores[$i]
}
